$wb = $excel.ActiveWorkbook

$oldGuid = "05e3eb0f-55c7-48b4-8b4f-120ac4583c02"
$newGuid = "185f8477-1d82-457c-9cf6-d4b222813430"
$addGuid = "f06f1f58-d76a-490e-9da1-1e75b1249941"

$dateMain    = "2016-08-25 03:01:22"
$dateZhStart = "2016-08-25 03:01:16"
$dateZhEnd   = "2016-08-25 03:01:34"
$dateDeEnd   = "2016-08-25 03:01:42"

$zhHashNew = "ecad8d5a1549ec58d999efce203a003077319e30"
$deHashNew = "ecad8d5a1549ec58d999efce203a003077319e30"
$zhHashAdd = "f39d8bd11bc593e915d4581b95a250a299905da0"
$deHashAdd = "f39d8bd11bc593e915d4581b95a250a299905da0"

$numFmtDate = "yyyy-mm-dd HH:mm:ss"

# =====================================================================
# Sheet "Overview"
# =====================================================================
$ov = $wb.Worksheets.Item("Overview")

# -- update row 2 (rename GUID, refresh date) --
$ov.Range("A2").Value = "$newGuid.md"
$ov.Range("B2").Value = "e2e\$newGuid.md"
$ov.Range("G2").Value = $dateMain

foreach ($hl in $ov.Hyperlinks) {
    $r = $hl.Range.Address()
    if ($r -eq "`$B`$2") {
        $hl.TextToDisplay = "e2e\$newGuid.md"
    }
}

# -- add row 3 (new file) --
$ov.Range("A3").Value = "$addGuid.md"
$ov.Range("B3").Value = "e2e\$addGuid.md"
$ov.Range("B3").Style = "Hyperlink"
$ov.Range("C3").Value = ".md"
$ov.Range("E3").Value = "Handed back: in sync with en-US"
$ov.Range("F3").Value = "Handed back: in sync with en-US"
$ov.Range("G3").Value = $dateMain
$ov.Range("G3").NumberFormat = $numFmtDate

$ov.Hyperlinks.Add(
    $ov.Range("B3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8cf5a67b44156cb17e05918bb4c76e9a10af2452/e2e/$addGuid.md",
    "",
    "",
    "e2e\$addGuid.md"
) | Out-Null

# =====================================================================
# Sheet "zh-cn"
# =====================================================================
$zh = $wb.Worksheets.Item("zh-cn")

# -- update row 2 --
$zh.Range("A2").Value = "$newGuid.md"
$zh.Range("G2").Value = "$newGuid.$zhHashNew.zh-cn.xlf"
$zh.Range("H2").Value = $dateZhStart
$zh.Range("I2").Value = "$newGuid.md"
$zh.Range("J2").Value = "$newGuid.$zhHashNew.zh-cn.xlf"
$zh.Range("K2").Value = $dateZhEnd

foreach ($hl in $zh.Hyperlinks) {
    $r = $hl.Range.Address()
    if ($r -eq "`$A`$2") {
        $hl.TextToDisplay = "$newGuid.md"
    } elseif ($r -eq "`$I`$2") {
        $hl.TextToDisplay = "$newGuid.md"
    }
}

# -- add row 3 --
$zh.Range("A3").Value = "$addGuid.md"
$zh.Range("A3").Style = "Hyperlink"
$zh.Range("B3").Value = ".md"
$zh.Range("C3").Value = "Handed back: in sync with en-US"
$zh.Range("D3").Value = "e2e"
$zh.Range("E3").Value = "ht"
$zh.Range("F3").Value = "'True"
$zh.Range("G3").Value = "$addGuid.$zhHashAdd.zh-cn.xlf"
$zh.Range("H3").Value = $dateZhStart
$zh.Range("H3").NumberFormat = $numFmtDate
$zh.Range("I3").Value = "$addGuid.md"
$zh.Range("I3").Style = "Hyperlink"
$zh.Range("J3").Value = "$addGuid.$zhHashAdd.zh-cn.xlf"
$zh.Range("K3").Value = $dateZhEnd
$zh.Range("K3").NumberFormat = $numFmtDate
$zh.Range("L3").Value = "'"
$zh.Range("M3").Value = "'True"
$zh.Range("N3").Value = "'"
$zh.Range("O3").Value = "'False"
$zh.Range("P3").Value = "'"

$zh.Hyperlinks.Add(
    $zh.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8cf5a67b44156cb17e05918bb4c76e9a10af2452/e2e/$addGuid.md",
    "",
    "",
    "$addGuid.md"
) | Out-Null

$zh.Hyperlinks.Add(
    $zh.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/fd998abfb3466da2a1b2d13092d8dd0688e6d2bd/e2e/$addGuid.md",
    "",
    "",
    "$addGuid.md"
) | Out-Null

# =====================================================================
# Sheet "de-de"
# =====================================================================
$de = $wb.Worksheets.Item("de-de")

# -- update row 2 --
$de.Range("A2").Value = "$newGuid.md"
$de.Range("G2").Value = "$newGuid.$deHashNew.de-de.xlf"
$de.Range("H2").Value = $dateMain
$de.Range("I2").Value = "$newGuid.md"
$de.Range("J2").Value = "$newGuid.$deHashNew.de-de.xlf"
$de.Range("K2").Value = $dateDeEnd

foreach ($hl in $de.Hyperlinks) {
    $r = $hl.Range.Address()
    if ($r -eq "`$A`$2") {
        $hl.TextToDisplay = "$newGuid.md"
    } elseif ($r -eq "`$I`$2") {
        $hl.TextToDisplay = "$newGuid.md"
    }
}

# -- add row 3 --
$de.Range("A3").Value = "$addGuid.md"
$de.Range("A3").Style = "Hyperlink"
$de.Range("B3").Value = ".md"
$de.Range("C3").Value = "Handed back: in sync with en-US"
$de.Range("D3").Value = "e2e"
$de.Range("E3").Value = "ht"
$de.Range("F3").Value = "'True"
$de.Range("G3").Value = "$addGuid.$deHashAdd.de-de.xlf"
$de.Range("H3").Value = $dateMain
$de.Range("H3").NumberFormat = $numFmtDate
$de.Range("I3").Value = "$addGuid.md"
$de.Range("I3").Style = "Hyperlink"
$de.Range("J3").Value = "$addGuid.$deHashAdd.de-de.xlf"
$de.Range("K3").Value = $dateDeEnd
$de.Range("K3").NumberFormat = $numFmtDate
$de.Range("L3").Value = "'"
$de.Range("M3").Value = "'True"
$de.Range("N3").Value = "'"
$de.Range("O3").Value = "'False"
$de.Range("P3").Value = "'"

$de.Hyperlinks.Add(
    $de.Range("A3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8cf5a67b44156cb17e05918bb4c76e9a10af2452/e2e/$addGuid.md",
    "",
    "",
    "$addGuid.md"
) | Out-Null

$de.Hyperlinks.Add(
    $de.Range("I3"),
    "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/4340b949579dfa81fd47284b1ab697dbd180551a/e2e/$addGuid.md",
    "",
    "",
    "$addGuid.md"
) | Out-Null

# =====================================================================
# Extend table ranges to include the new row
# =====================================================================
$ov.ListObjects.Item(1).Resize($ov.Range("A1:G3"))
$zh.ListObjects.Item(1).Resize($zh.Range("A1:P3"))
$de.ListObjects.Item(1).Resize($de.Range("A1:P3"))

Write-Host "Edit complete"
